# Apply updated dSF (column F) values on Sheet1, per repull/mean-calculation update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row number -> new value for column F (dSF)
$updates = @{
    2  = 1
    4  = -2
    6  = 5
    7  = 5
    8  = -2
    9  = -2
    10 = -1
    11 = 3
    12 = -4
    14 = 1
    15 = -3
    16 = 1
    17 = -3
    20 = -3
    21 = 9
    22 = 4
    23 = -3
    24 = -3
    25 = 2
    26 = 3
    27 = -3
    28 = -1
    29 = -1
    30 = -2
    31 = -1
    33 = -3
    34 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
